$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the team record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style from an existing header cell (e.g. AC1) so the new
# headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 76   # AD
    $ws.Cells.Item($row, 31).Value = 86   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
